$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 796, shifting existing rows 796:837 down to 797:838
$ws.Rows.Item(796).Insert(-4121)  # -4121 = xlShiftDown

# Populate the newly inserted row with the new data point.
# Prefix the date-looking value with an apostrophe so it is stored as plain
# text (not auto-converted to a date serial number), then reset the cell
# style back to Normal so no stray number-format style is left behind.
$ws.Cells.Item(796, 1).Value = "'2026/02/11"
$ws.Cells.Item(796, 1).Style = "Normal"
$ws.Cells.Item(796, 2).Value = "水"
$ws.Cells.Item(796, 3).Value = 0
$ws.Cells.Item(796, 4).Value = 201
